$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 / column B ("comfort_food") currently holds the text
# "frozen yogurt, pizza, fast food" (shared string). Replace it with the
# numeric value 5555 — this also makes that shared string unused, so
# Excel drops it from sharedStrings.xml and every later shared-string
# index shifts down by one automatically.
$ws.Range("B4").Value = 5555

# Move the active selection from F8 to B5.
$ws.Range("B5").Select() | Out-Null
